# Add I0 and IF columns to the worksheet (headers + data)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of the last existing header cell (H1) onto the
# new header cells so they match the other bold/bordered/centered headers.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# Header text
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for rows 2-42, columns I (I0) and J (IF)
$arr = New-Object 'object[,]' 41,2
$arr[0,0] = 9
$arr[0,1] = 9
$arr[1,0] = 8
$arr[1,1] = 8
$arr[2,0] = 8
$arr[2,1] = 8
$arr[3,0] = 8
$arr[3,1] = 8
$arr[4,0] = 6
$arr[4,1] = 6
$arr[5,0] = 8
$arr[5,1] = 8
$arr[6,0] = 8
$arr[6,1] = 8
$arr[7,0] = 9
$arr[7,1] = 9
$arr[8,0] = 9
$arr[8,1] = 9
$arr[9,0] = 9
$arr[9,1] = 9
$arr[10,0] = 9
$arr[10,1] = 9
$arr[11,0] = 9
$arr[11,1] = 9
$arr[12,0] = 9
$arr[12,1] = 9
$arr[13,0] = 11
$arr[13,1] = 11
$arr[14,0] = 10
$arr[14,1] = 10
$arr[15,0] = 9
$arr[15,1] = 9
$arr[16,0] = 9
$arr[16,1] = 9
$arr[17,0] = 9
$arr[17,1] = 9
$arr[18,0] = 9
$arr[18,1] = 9
$arr[19,0] = 9
$arr[19,1] = 9
$arr[20,0] = 9
$arr[20,1] = 9
$arr[21,0] = 9
$arr[21,1] = 9
$arr[22,0] = 9
$arr[22,1] = 9
$arr[23,0] = 9
$arr[23,1] = 9
$arr[24,0] = 9
$arr[24,1] = 9
$arr[25,0] = 9
$arr[25,1] = 9
$arr[26,0] = 9
$arr[26,1] = 9
$arr[27,0] = 10
$arr[27,1] = 10
$arr[28,0] = 9
$arr[28,1] = 9
$arr[29,0] = 9
$arr[29,1] = 9
$arr[30,0] = 9
$arr[30,1] = 9
$arr[31,0] = 8
$arr[31,1] = 9
$arr[32,0] = 9
$arr[32,1] = 9
$arr[33,0] = 9
$arr[33,1] = 9
$arr[34,0] = 9
$arr[34,1] = 9
$arr[35,0] = 8
$arr[35,1] = 8
$arr[36,0] = 7
$arr[36,1] = 7
$arr[37,0] = 9
$arr[37,1] = 9
$arr[38,0] = 6
$arr[38,1] = 6
$arr[39,0] = 6
$arr[39,1] = 6
$arr[40,0] = 3
$arr[40,1] = 3

$ws.Range("I2:J42").Value = $arr

$ws.Range("A1").Select() | Out-Null
